# Fix: source delimiter issues
#
# The "data" sheet's table (Table1) had a "sources" column (G) that always
# held the literal multiplier 1.2 with no indication of whether a given
# data row's title (B) refers to an actual country or to an aggregate
# region (e.g. "World", "European Union"). This adds a new helper column F
# that labels each row as "Country" or "Region" (the delimiter used when
# building the source citation text), and marks the existing "sources"
# column (G) as text-formatted.
#
# It also restores the "data" sheet as the active/selected tab (it had
# drifted to "info"), with the selection left on D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$firstDataRow = 2
$lastDataRow = 152

# Rows whose title (column B) is an aggregate region rather than a single
# country - these get the "Region" label instead of "Country".
$regionRows = @(148, 152)

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    if ($regionRows -contains $r) {
        $ws.Range("F$r").Value = "Region"
    } else {
        $ws.Range("F$r").Value = "Country"
    }
}

# The "sources" column (G) keeps its existing values, but should be
# formatted as text.
$ws.Range("G$firstDataRow`:G$lastDataRow").NumberFormat = "@"

# Re-select the "data" sheet as the active tab (previously "info" was left
# selected) and move the active cell/selection to D11.
[void]$ws.Activate()
$ws.Range("D11").Select() | Out-Null

# The "info" sheet should no longer be the tab-selected sheet.
$infoWs = $wb.Worksheets.Item("info")
$infoWs.Range("B2").Select() | Out-Null
[void]$ws.Activate()
